$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "58.375.36"
    "E2" = "  -2.68%  "
    "D3" = "2.453.18"
    "E3" = "  -3.71%  "
    "D4" = "1.00"
    "E4" = "  -0.10%  "
    "D5" = "528.07"
    "E5" = "  -1.93%  "
    "D6" = "133.87"
    "E6" = "  -7.00%  "
    "D7" = "0.998"
    "E7" = "  +0.18%  "
    "D8" = "0.557"
    "E8" = "  -2.42%  "
    "D9" = "2.458.32"
    "E9" = "  -4.18%  "
    "E10" = "  -2.35%  "
    "E11" = "  -0.33%  "
    "D12" = "5.33"
    "E12" = "  -2.77%  "
    "E13" = "  -5.45%  "
    "D14" = "2.887.59"
    "E14" = "  -3.73%  "
    "D15" = "58.299.69"
    "E15" = "  -2.77%  "
    "D16" = "22.63"
    "E16" = "  -5.88%  "
    "D17" = "0.0000138"
    "E17" = "  -4.01%  "
    "D18" = "2.460.64"
    "E18" = "  -4.59%  "
    "D19" = "10.71"
    "E19" = "  -4.73%  "
    "B20" = "BitcoinCash"
    "C20" = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
    "D20" = "320.95"
    "E20" = "  -1.75%  "
    "B21" = "Polkadot"
    "C21" = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
    "D21" = "4.18"
    "E21" = "  -3.40%  "
    "E22" = "  -0.34%  "
    "D23" = "5.71"
    "E23" = "  -4.29%  "
    "D24" = "62.56"
    "E24" = "  -1.15%  "
    "E25" = "  -6.06%  "
    "E26" = "  -1.19%  "
    "D27" = "0.983"
    "E27" = "  -1.33%  "
    "D28" = "7.42"
    "E28" = "  -7.53%  "
    "D29" = "0.0₃0748"
    "E29" = "  -5.84%  "
    "D30" = "6.47"
    "E30" = "  -8.13%  "
    "E31" = "  -3.53%  "
    "D32" = "163.89"
    "E32" = "  -0.68%  "
    "E33" = "  +0.06%  "
    "D34" = "1.10"
    "E34" = "  -7.56%  "
    "B35" = "EthereumClassic"
    "C35" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    "D35" = "18.22"
    "E35" = "  -2.63%  "
    "B36" = "ImmutableX"
    "C36" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "D36" = "1.35"
    "E36" = "  -8.77%  "
    "D37" = "4.00"
    "E37" = "  -8.82%  "
    "D38" = "1.53"
    "E38" = "  -6.06%  "
    "D39" = "36.42"
    "E39" = "  -1.63%  "
    "D40" = "0.804"
    "E40" = "  -3.62%  "
    "D41" = "3.53"
    "E41" = "  -5.07%  "
    "D42" = "274.37"
    "E42" = "  -9.22%  "
    "E43" = "  +0.12%  "
    "D44" = "5.06"
    "E44" = "  -9.24%  "
    "E45" = "  -0.06%  "
    "E46" = "  -4.08%  "
    "D47" = "0.0922"
    "E47" = "  -1.67%  "
    "D48" = "120.06"
    "E48" = "  -5.59%  "
    "D49" = "0.0505"
    "E49" = "  -2.66%  "
    "D50" = "0.0218"
    "E50" = "  -4.54%  "
    "D51" = "17.00"
    "E51" = "  -6.55%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = $origStyle
}